$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "---"
$ws.Range("D4").Value = "Service providers performance: ---"

$ws.Range("B5").Value = "---"
$ws.Range("D5").Value = "Investment status: ---"

$ws.Range("B6").Value = "---"
$ws.Range("D6").Value = "Lessons Learned: ---"
